$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct value assignments for cells whose new text is unambiguous (not a plain number)
$ws.Range("D2").Value = "69.879.19"
$ws.Range("D3").Value = "2.565.94"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "2.564.76"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  +11.53%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("E14").Value = "  +5.79%  "
$ws.Range("D15").Value = "3.005.15"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "69.750.94"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("D18").Value = "2.573.60"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  +3.25%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "0.0₃0917"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("E38").Value = "  +3.09%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  +1.54%  "

# Helper to force a numeric-looking string to be stored as text, matching the
# original cells which are plain text (not real numbers), while restoring the
# cell style back to the sheet default afterwards so no visible formatting changes.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Cells whose new value looks like a plain decimal number; without forcing text
# Excel would convert them into numeric cells, which would not match the source data.
Set-TextValue "D5" "600.61"
Set-TextValue "D6" "178.13"
Set-TextValue "D8" "0.519"
Set-TextValue "D12" "0.344"
Set-TextValue "D13" "5.01"
Set-TextValue "D14" "0.0000182"
Set-TextValue "D16" "26.27"
Set-TextValue "D20" "11.18"
Set-TextValue "D21" "365.54"
Set-TextValue "D24" "70.81"
Set-TextValue "D31" "513.37"
Set-TextValue "D32" "7.80"
Set-TextValue "D36" "164.06"
Set-TextValue "D38" "18.99"
Set-TextValue "D40" "1.35"
Set-TextValue "D41" "1.76"
Set-TextValue "D43" "4.92"
Set-TextValue "D45" "2.46"
Set-TextValue "D46" "39.03"
Set-TextValue "D47" "151.84"
Set-TextValue "D48" "3.62"
